$wb = $excel.ActiveWorkbook
$basic = $wb.Worksheets.Item("Basic")

# Fill in the previously-empty D:H columns for rows 10-18 on the "Basic" sheet
$rows = @(
    @(10, 11, 309501013, 100000096, 80098733, 0.04888888888888889),
    @(11, 14, 330785622, 100000030, 78511175, 0.20518518518518516),
    @(12, 14, 335953465, 100000040, 79041186, 0.21305555555555555),
    @(13, 13, 341161320, 100000068, 79339459, 0.2441087962962963),
    @(14, 13, 341846423, 100000073, 79774012, 0.2865740740740741),
    @(15, 13, 386115817, 100000045, 80671850, 0.5291319444444444),
    @(16, 13, 367708286, 100000121, 82091596, 0.5789930555555556),
    @(17, 11, 395682063, 100000075, 81634279, 0.7478472222222222),
    @(18, 12, 401141644, 100000047, 82182689, 0.7575115740740741)
)

$hFormat = $basic.Range("H3").NumberFormat

foreach ($item in $rows) {
    $r = $item[0]
    $basic.Range("D$r").Value = $item[1]
    $basic.Range("E$r").Value = $item[2]
    $basic.Range("F$r").Value = $item[3]
    $basic.Range("G$r").Value = $item[4]
    $basic.Range("H$r").Value = $item[5]
    $basic.Range("H$r").NumberFormat = $hFormat
}

# Move selection on the "Basic" sheet
$basic.Activate()
$basic.Range("G17").Select()

# Add the new "Durable" sheet after the last sheet, and make it the active tab
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$durable = $wb.Worksheets.Add($null, $lastSheet)
$durable.Name = "Durable"
$durable.Activate()
$durable.Range("H32").Select()
